$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is inserted as row 30; the previously-existing
# rows 30-33 shift down to 31-34 (their contents stay the same, only their
# row numbers change - which Excel does automatically on Insert).
$ws.Rows.Item(30).Insert()

# Fill in the newly-inserted row 30 with the new week's data.
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44644
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100114007
$ws.Range("G30").Value = "Jengibre"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15500
$ws.Range("N30").Value = "$/caja 13 kilos"
$ws.Range("O30").Value = "Perú"
$ws.Range("P30").Value = 1192
$ws.Range("Q30").Value = 13
$ws.Range("R30").Value = "Hortaliza"
